$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 (R) values ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 197
$wsOff.Range("C3").Value = 122
$wsOff.Range("D3").Value = 62
$wsOff.Range("E3").Value = 25
$wsOff.Range("F3").Value = 5
$wsOff.Range("G3").Value = 3

# --- DEF sheet: row 3 (R) values ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 217
$wsDef.Range("C3").Value = 146
$wsDef.Range("D3").Value = 43
$wsDef.Range("F3").Value = 7
